# Apply updated market-price figures (currentAveragePrice / LevePrice / LeveProfit columns)
# pulled from the scheduled market-data refresh, per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 994.0909
$ws.Range("I19").Value = 962.6667
$ws.Range("J19").Value = 1005.875
$ws.Range("K19").Value = 962.6667
$ws.Range("L19").Value = 1005.875
$ws.Range("M19").Value = -787.6667
$ws.Range("N19").Value = -1355.875
$ws.Range("H43").Value = 3851653.5
$ws.Range("I43").Value = 5132204.5
$ws.Range("K43").Value = 5132204.5
$ws.Range("M43").Value = -5132135.5
$ws.Range("H138").Value = 4787.409
$ws.Range("I138").Value = 1435.2222
$ws.Range("J138").Value = 5316.7017
$ws.Range("K138").Value = 4305.6666
$ws.Range("L138").Value = 15950.1051
$ws.Range("M138").Value = 834.3334000000004
$ws.Range("N138").Value = -26230.1051

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 18000
$ws.Range("I37").Value = 18000
$ws.Range("K37").Value = 18000
$ws.Range("M37").Value = -17727
$ws.Range("H61").Value = 20996.875
$ws.Range("I61").Value = 20996.875
$ws.Range("K61").Value = 20996.875
$ws.Range("M61").Value = -20784.875
$ws.Range("H74").Value = 8335006
$ws.Range("I74").Value = 8622282
$ws.Range("K74").Value = 8622282
$ws.Range("M74").Value = -8621408
$ws.Range("H77").Value = 8335006
$ws.Range("I77").Value = 8622282
$ws.Range("K77").Value = 43111410
$ws.Range("M77").Value = -43107042
$ws.Range("H122").Value = 4112.8823
$ws.Range("I122").Value = 2502.375
$ws.Range("J122").Value = 5544.4443
$ws.Range("K122").Value = 7507.125
$ws.Range("L122").Value = 16633.3329
$ws.Range("M122").Value = -5057.125
$ws.Range("N122").Value = -21533.3329
$ws.Range("H132").Value = 4780.892
$ws.Range("I132").Value = 1599.7142
$ws.Range("K132").Value = 4799.142599999999
$ws.Range("M132").Value = -2269.142599999999
$ws.Range("H136").Value = 20996.875
$ws.Range("I136").Value = 20996.875
$ws.Range("K136").Value = 62990.625
$ws.Range("M136").Value = -60440.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1308.5652
$ws.Range("I86").Value = 1339
$ws.Range("J86").Value = 1239
$ws.Range("K86").Value = 1339
$ws.Range("L86").Value = 1239
$ws.Range("M86").Value = -216
$ws.Range("N86").Value = -3485
$ws.Range("H89").Value = 1308.5652
$ws.Range("I89").Value = 1339
$ws.Range("J89").Value = 1239
$ws.Range("K89").Value = 6695
$ws.Range("L89").Value = 6195
$ws.Range("M89").Value = -1079
$ws.Range("N89").Value = -17427
$ws.Range("H94").Value = 857182.3
$ws.Range("I94").Value = 1054624.4
$ws.Range("J94").Value = 1600
$ws.Range("K94").Value = 1054624.4
$ws.Range("L94").Value = 1600
$ws.Range("M94").Value = -1054173.4
$ws.Range("N94").Value = -2502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 667.875
$ws.Range("I22").Value = 732.5
$ws.Range("J22").Value = 629.1
$ws.Range("K22").Value = 732.5
$ws.Range("L22").Value = 629.1
$ws.Range("M22").Value = -382.5
$ws.Range("N22").Value = -1329.1
$ws.Range("H58").Value = 557190.3
$ws.Range("I58").Value = 909878.9399999999
$ws.Range("J58").Value = 2965.4285
$ws.Range("K58").Value = 909878.9399999999
$ws.Range("L58").Value = 2965.4285
$ws.Range("M58").Value = -909675.9399999999
$ws.Range("N58").Value = -3371.4285
$ws.Range("H136").Value = 557190.3
$ws.Range("I136").Value = 909878.9399999999
$ws.Range("J136").Value = 2965.4285
$ws.Range("K136").Value = 2729636.82
$ws.Range("L136").Value = 8896.2855
$ws.Range("M136").Value = -2727086.82
$ws.Range("N136").Value = -13996.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 266.2
$ws.Range("J97").Value = 200
$ws.Range("L97").Value = 600
$ws.Range("N97").Value = -1592
$ws.Range("H132").Value = 9222.857
$ws.Range("J132").Value = 15327.75
$ws.Range("L132").Value = 137949.75
$ws.Range("N132").Value = -143009.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3407103
$ws.Range("J70").Value = 6999.8
$ws.Range("L70").Value = 6999.8
$ws.Range("N70").Value = -7539.8
$ws.Range("H73").Value = 3407103
$ws.Range("J73").Value = 6999.8
$ws.Range("L73").Value = 6999.8
$ws.Range("N73").Value = -8871.799999999999
$ws.Range("H102").Value = 29422500
$ws.Range("I102").Value = 45466584
$ws.Range("J102").Value = 8349.833000000001
$ws.Range("K102").Value = 45466584
$ws.Range("L102").Value = 8349.833000000001
$ws.Range("M102").Value = -45464962
$ws.Range("N102").Value = -11593.833
$ws.Range("H132").Value = 187662.73
$ws.Range("I132").Value = 291353.84
$ws.Range("K132").Value = 874061.52
$ws.Range("M132").Value = -871531.52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1120.4814
$ws.Range("I22").Value = 991.1429000000001
$ws.Range("J22").Value = 1259.7693
$ws.Range("K22").Value = 991.1429000000001
$ws.Range("L22").Value = 1259.7693
$ws.Range("M22").Value = -696.1429000000001
$ws.Range("N22").Value = -1849.7693
$ws.Range("H27").Value = 1120.4814
$ws.Range("I27").Value = 991.1429000000001
$ws.Range("J27").Value = 1259.7693
$ws.Range("K27").Value = 991.1429000000001
$ws.Range("L27").Value = 1259.7693
$ws.Range("M27").Value = -884.1429000000001
$ws.Range("N27").Value = -1473.7693
$ws.Range("H46").Value = 6267.0713
$ws.Range("J46").Value = 6364.5386
$ws.Range("L46").Value = 6364.5386
$ws.Range("N46").Value = -6740.5386
$ws.Range("H55").Value = 441.625
$ws.Range("J55").Value = 551.1667
$ws.Range("L55").Value = 551.1667
$ws.Range("N55").Value = -897.1667
$ws.Range("H93").Value = 1650.4445
$ws.Range("I93").Value = 1229.8334
$ws.Range("J93").Value = 2491.6667
$ws.Range("K93").Value = 1229.8334
$ws.Range("L93").Value = 2491.6667
$ws.Range("M93").Value = 18.16660000000002
$ws.Range("N93").Value = -4987.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2987582.2
$ws.Range("I81").Value = 4168026.5
$ws.Range("J81").Value = 2331780
$ws.Range("K81").Value = 8336053
$ws.Range("L81").Value = 4663560
$ws.Range("M81").Value = -8334992
$ws.Range("N81").Value = -4665682
$ws.Range("H84").Value = 2987582.2
$ws.Range("I84").Value = 4168026.5
$ws.Range("J84").Value = 2331780
$ws.Range("K84").Value = 41680265
$ws.Range("L84").Value = 23317800
$ws.Range("M84").Value = -41674961
$ws.Range("N84").Value = -23328408
$ws.Range("H132").Value = 19161802
$ws.Range("I132").Value = 3087199.5
$ws.Range("K132").Value = 9261598.5
$ws.Range("M132").Value = -9259068.5
